# table_1_78.xlsx - "progress is being made"
# Change the per-row divisor formula in column D (rows 35-38) from the
# shared formula "=$D$3/3" to an individual formula "=D3/2", which
# ripples through the dependent E/K columns and the C40 grand total via
# normal recalculation. Also update the active window's scroll position
# and selection to match the author's view state when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D35").Formula = "=D3/2"
$ws.Range("D36").Formula = "=D3/2"
$ws.Range("D37").Formula = "=D3/2"
$ws.Range("D38").Formula = "=D3/2"

# Move the viewport so row 10 is at the top-left and select D37, matching
# the saved sheet view.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D37").Select()
